$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 531.8
$ws.Range("I9").Value = 605.4286
$ws.Range("J9").Value = 467.375
$ws.Range("K9").Value = 605.4286
$ws.Range("L9").Value = 467.375
$ws.Range("M9").Value = -436.4286
$ws.Range("N9").Value = -805.375

$ws.Range("H17").Value = 1735.75
$ws.Range("J17").Value = 2148.6155
$ws.Range("L17").Value = 6445.8465
$ws.Range("N17").Value = -6781.8465

$ws.Range("H52").Value = 522.5
$ws.Range("I52").Value = 522.5
$ws.Range("K52").Value = 1567.5
$ws.Range("M52").Value = -1407.5

$ws.Range("H125").Value = 1614.5834
$ws.Range("I125").Value = 1945.8334
$ws.Range("J125").Value = 1283.3334
$ws.Range("K125").Value = 17512.5006
$ws.Range("L125").Value = 11550.0006
$ws.Range("M125").Value = -15052.5006
$ws.Range("N125").Value = -16470.0006

$ws.Range("H127").Value = 647
$ws.Range("I127").Value = 647
$ws.Range("K127").Value = 1941
$ws.Range("M127").Value = 3019

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H130").Value = 130000
$ws.Range("J130").Value = 130000
$ws.Range("L130").Value = 130000
$ws.Range("N130").Value = -140040

$ws.Range("H135").Value = 1150.7941
$ws.Range("I135").Value = 1250.1333
$ws.Range("K135").Value = 11251.1997
$ws.Range("M135").Value = -8716.199699999999

$ws.Range("H137").Value = 2210.15
$ws.Range("I137").Value = 1850.4667
$ws.Range("J137").Value = 3289.2
$ws.Range("K137").Value = 5551.4001
$ws.Range("L137").Value = 9867.599999999999
$ws.Range("M137").Value = -3001.4001
$ws.Range("N137").Value = -14967.6

$ws.Range("H138").Value = 1642.2858
$ws.Range("J138").Value = 3604.7778
$ws.Range("L138").Value = 10814.3334
$ws.Range("N138").Value = -21094.3334

$ws.Range("H141").Value = 1397.5
$ws.Range("I141").Value = 1397.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4192.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 987.5
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 55562244
$ws.Range("I32").Value = 62505024
$ws.Range("K32").Value = 62505024
$ws.Range("M32").Value = -62504737

$ws.Range("H45").Value = 4066.9
$ws.Range("I45").Value = 4318.8887
$ws.Range("K45").Value = 4318.8887
$ws.Range("M45").Value = -3941.8887

$ws.Range("H132").Value = 2906
$ws.Range("I132").Value = 2893.0784
$ws.Range("J132").Value = 3125.6667
$ws.Range("K132").Value = 8679.235199999999
$ws.Range("L132").Value = 9377.000100000001
$ws.Range("M132").Value = -6149.235199999999
$ws.Range("N132").Value = -14437.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3112.4285
$ws.Range("I86").Value = 3157.4
$ws.Range("K86").Value = 3157.4
$ws.Range("M86").Value = -2034.4

$ws.Range("H89").Value = 3112.4285
$ws.Range("I89").Value = 3157.4
$ws.Range("K89").Value = 15787
$ws.Range("M89").Value = -10171

$ws.Range("H94").Value = 1629.3077
$ws.Range("J94").Value = 1530.4546
$ws.Range("L94").Value = 1530.4546
$ws.Range("N94").Value = -2432.4546

$ws.Range("H112").Value = 90000
$ws.Range("J112").Value = 90000
$ws.Range("L112").Value = 90000
$ws.Range("N112").Value = -92954

$ws.Range("H134").Value = 2923.5144
$ws.Range("I134").Value = 1495.36
$ws.Range("J134").Value = 6493.9
$ws.Range("K134").Value = 4486.08
$ws.Range("L134").Value = 19481.7
$ws.Range("M134").Value = -1951.08
$ws.Range("N134").Value = -24551.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3393
$ws.Range("I22").Value = 491.8
$ws.Range("K22").Value = 491.8
$ws.Range("M22").Value = -141.8

$ws.Range("H31").Value = 81412.586
$ws.Range("I31").Value = 112638.97
$ws.Range("K31").Value = 112638.97
$ws.Range("M31").Value = -112343.97

$ws.Range("H34").Value = 81412.586
$ws.Range("I34").Value = 112638.97
$ws.Range("K34").Value = 112638.97
$ws.Range("M34").Value = -112436.97

$ws.Range("H58").Value = 3730.6155
$ws.Range("I58").Value = 1179.1666
$ws.Range("K58").Value = 1179.1666
$ws.Range("M58").Value = -976.1666

$ws.Range("H99").Value = 8465.157999999999
$ws.Range("I99").Value = 3453.7273
$ws.Range("J99").Value = 15355.875
$ws.Range("K99").Value = 3453.7273
$ws.Range("L99").Value = 15355.875
$ws.Range("M99").Value = -1955.7273
$ws.Range("N99").Value = -18351.875

$ws.Range("H126").Value = 8465.157999999999
$ws.Range("I126").Value = 3453.7273
$ws.Range("J126").Value = 15355.875
$ws.Range("K126").Value = 10361.1819
$ws.Range("L126").Value = 46067.625
$ws.Range("M126").Value = -7891.1819
$ws.Range("N126").Value = -51007.625

$ws.Range("H132").Value = 5685266.5
$ws.Range("I132").Value = 3554.561
$ws.Range("K132").Value = 10663.683
$ws.Range("M132").Value = -8133.683000000001

$ws.Range("H134").Value = 9543.244000000001
$ws.Range("I134").Value = 11248.906
$ws.Range("K134").Value = 33746.718
$ws.Range("M134").Value = -31211.718

$ws.Range("H136").Value = 3730.6155
$ws.Range("I136").Value = 1179.1666
$ws.Range("K136").Value = 3537.4998
$ws.Range("M136").Value = -987.4998000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7727.8237
$ws.Range("I56").Value = 7727.8237
$ws.Range("K56").Value = 7727.8237
$ws.Range("M56").Value = -7197.8237

$ws.Range("H103").Value = 1147.8572
$ws.Range("I103").Value = 52
$ws.Range("J103").Value = 2609
$ws.Range("K103").Value = 156
$ws.Range("L103").Value = 7827
$ws.Range("M103").Value = 723
$ws.Range("N103").Value = -9585

$ws.Range("H131").Value = 13019.9
$ws.Range("J131").Value = 18248.715
$ws.Range("L131").Value = 54746.145
$ws.Range("N131").Value = -64826.145

$ws.Range("H133").Value = 6709.357
$ws.Range("I133").Value = 8482.75
$ws.Range("K133").Value = 25448.25
$ws.Range("M133").Value = -20388.25

$ws.Range("H134").Value = 1000
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = 2070

$ws.Range("H140").Value = 2712.0789
$ws.Range("I140").Value = 2168.6667
$ws.Range("J140").Value = 4749.875
$ws.Range("K140").Value = 6506.000100000001
$ws.Range("L140").Value = 14249.625
$ws.Range("M140").Value = -1326.000100000001
$ws.Range("N140").Value = -24609.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 128.96774
$ws.Range("I2").Value = 73.23529000000001
$ws.Range("K2").Value = 73.23529000000001
$ws.Range("M2").Value = 39.76470999999999

$ws.Range("H10").Value = 1000000000
$ws.Range("I10").Value = 1000000000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1000000000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -999999831
$ws.Range("N10").ClearContents()

$ws.Range("H132").Value = 403088.22
$ws.Range("I132").Value = 610289.9
$ws.Range("J132").Value = 26357.908
$ws.Range("K132").Value = 1830869.7
$ws.Range("L132").Value = 79073.724
$ws.Range("M132").Value = -1828339.7
$ws.Range("N132").Value = -84133.724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2863.5151
$ws.Range("I40").Value = 2607.44
$ws.Range("J40").Value = 3663.75
$ws.Range("K40").Value = 2607.44
$ws.Range("L40").Value = 3663.75
$ws.Range("M40").Value = -2471.44
$ws.Range("N40").Value = -3935.75

$ws.Range("H68").Value = 3812.647
$ws.Range("I68").Value = 2042.8572
$ws.Range("J68").Value = 5051.5
$ws.Range("K68").Value = 2042.8572
$ws.Range("L68").Value = 5051.5
$ws.Range("M68").Value = -1293.8572
$ws.Range("N68").Value = -6549.5

$ws.Range("H71").Value = 3812.647
$ws.Range("I71").Value = 2042.8572
$ws.Range("J71").Value = 5051.5
$ws.Range("K71").Value = 10214.286
$ws.Range("L71").Value = 25257.5
$ws.Range("M71").Value = -6470.286
$ws.Range("N71").Value = -32745.5

$ws.Range("H82").Value = 1572.909
$ws.Range("I82").Value = 1086
$ws.Range("K82").Value = 1086
$ws.Range("M82").Value = -725

$ws.Range("H85").Value = 1572.909
$ws.Range("I85").Value = 1086
$ws.Range("K85").Value = 1086
$ws.Range("M85").Value = 162

$ws.Range("H102").Value = 90561
$ws.Range("J102").Value = 90561
$ws.Range("L102").Value = 90561
$ws.Range("N102").Value = -97051

$ws.Range("H136").Value = 43658.633
$ws.Range("I136").Value = 2726.963
$ws.Range("K136").Value = 8180.889000000001
$ws.Range("M136").Value = -5630.889000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1389.6666
$ws.Range("I81").Value = 1389.6666
$ws.Range("K81").Value = 2779.3332
$ws.Range("M81").Value = -1718.3332

$ws.Range("H84").Value = 1389.6666
$ws.Range("I84").Value = 1389.6666
$ws.Range("K84").Value = 13896.666
$ws.Range("M84").Value = -8592.666000000001

$ws.Range("H132").Value = 511828.53
$ws.Range("I132").Value = 708123.5600000001
$ws.Range("K132").Value = 2124370.68
$ws.Range("M132").Value = -2121840.68

$ws.Range("H136").Value = 8846177
$ws.Range("I136").Value = 10270842
$ws.Range("K136").Value = 30812526
$ws.Range("M136").Value = -30809976
